$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1440.75
$ws.Range("I43").Value = 1635.7142
$ws.Range("J43").Value = 1289.1111
$ws.Range("K43").Value = 1635.7142
$ws.Range("L43").Value = 1289.1111
$ws.Range("M43").Value = -1566.7142
$ws.Range("N43").Value = -1427.1111

$ws.Range("H121").Value = 1300
$ws.Range("I121").Value = 633.3333
$ws.Range("J121").Value = 1433.3334
$ws.Range("K121").Value = 1899.9999
$ws.Range("L121").Value = 4300.0002
$ws.Range("M121").Value = -152.9999
$ws.Range("N121").Value = -7794.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17860798
$ws.Range("I32").Value = 20411788
$ws.Range("K32").Value = 20411788
$ws.Range("M32").Value = -20411501

$ws.Range("H74").Value = 5210.6
$ws.Range("I74").Value = 846.4
$ws.Range("J74").Value = 11756.9
$ws.Range("K74").Value = 846.4
$ws.Range("L74").Value = 11756.9
$ws.Range("M74").Value = 27.60000000000002
$ws.Range("N74").Value = -13504.9

$ws.Range("H77").Value = 5210.6
$ws.Range("I77").Value = 846.4
$ws.Range("J77").Value = 11756.9
$ws.Range("K77").Value = 4232
$ws.Range("L77").Value = 58784.5
$ws.Range("M77").Value = 136
$ws.Range("N77").Value = -67520.5

$ws.Range("H113").Value = 45000
$ws.Range("J113").Value = 45000
$ws.Range("L113").Value = 45000
$ws.Range("N113").Value = -53678

$ws.Range("H132").Value = 5170.8057
$ws.Range("I132").Value = 6745.6665
$ws.Range("J132").Value = 2966
$ws.Range("K132").Value = 20236.9995
$ws.Range("L132").Value = 8898
$ws.Range("M132").Value = -17706.9995
$ws.Range("N132").Value = -13958

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 716.2258
$ws.Range("I94").Value = 617.72
$ws.Range("J94").Value = 1126.6666
$ws.Range("K94").Value = 617.72
$ws.Range("L94").Value = 1126.6666
$ws.Range("M94").Value = -166.72
$ws.Range("N94").Value = -2028.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1021.4
$ws.Range("I35").Value = 1021.4
$ws.Range("K35").Value = 1021.4
$ws.Range("M35").Value = -727.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 5320
$ws.Range("J96").Value = 5320
$ws.Range("L96").Value = 15960
$ws.Range("N96").Value = -20078

$ws.Range("H129").Value = 1501.381
$ws.Range("J129").Value = 1796.25
$ws.Range("L129").Value = 5388.75
$ws.Range("N129").Value = -15388.75

$ws.Range("H131").Value = 572.92
$ws.Range("I131").Value = 270.4717
$ws.Range("J131").Value = 913.9787
$ws.Range("K131").Value = 811.4150999999999
$ws.Range("L131").Value = 2741.9361
$ws.Range("M131").Value = 4228.5849
$ws.Range("N131").Value = -12821.9361

$ws.Range("H134").Value = 62700.668
$ws.Range("I134").Value = 62700.668
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 188102.004
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -183032.004

$ws.Range("H139").Value = 480487.97
$ws.Range("I139").Value = 580536.3
$ws.Range("K139").Value = 1741608.9
$ws.Range("M139").Value = -1736468.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 26.166666
$ws.Range("J2").Value = 31.5
$ws.Range("L2").Value = 31.5
$ws.Range("N2").Value = -257.5

$ws.Range("H10").Value = 669266.7
$ws.Range("I10").Value = 1001900
$ws.Range("K10").Value = 1001900
$ws.Range("M10").Value = -1001731

$ws.Range("H11").Value = 70004
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 70004
$ws.Range("K11").Value = 0
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = 70004
$ws.Range("N11").Value = -70282

$ws.Range("H14").Value = 572076.9
$ws.Range("I14").Value = 572076.9
$ws.Range("K14").Value = 572076.9
$ws.Range("M14").Value = -571908.9

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0

$ws.Range("H102").Value = 3590.375
$ws.Range("I102").Value = 3090.182
$ws.Range("J102").Value = 4690.8
$ws.Range("K102").Value = 3090.182
$ws.Range("L102").Value = 4690.8
$ws.Range("M102").Value = -1468.182
$ws.Range("N102").Value = -7934.8

$ws.Range("H132").Value = 6915.385
$ws.Range("I132").Value = 7886.421
$ws.Range("J132").Value = 4279.7144
$ws.Range("K132").Value = 23659.263
$ws.Range("L132").Value = 12839.1432
$ws.Range("M132").Value = -21129.263
$ws.Range("N132").Value = -17899.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1376.5714
$ws.Range("I7").Value = 1251.5
$ws.Range("J7").Value = 1543.3334
$ws.Range("K7").Value = 1251.5
$ws.Range("L7").Value = 1543.3334
$ws.Range("M7").Value = -1139.5
$ws.Range("N7").Value = -1767.3334

$ws.Range("H94").Value = 25466
$ws.Range("J94").Value = 25466
$ws.Range("L94").Value = 25466
$ws.Range("N94").Value = -26818

$ws.Range("H126").Value = 1376.5714
$ws.Range("I126").Value = 1251.5
$ws.Range("J126").Value = 1543.3334
$ws.Range("K126").Value = 3754.5
$ws.Range("L126").Value = 4630.0002
$ws.Range("M126").Value = -1284.5
$ws.Range("N126").Value = -9570.0002

$ws.Range("H132").Value = 5861.5405
$ws.Range("I132").Value = 7127.9546
$ws.Range("J132").Value = 4004.1333
$ws.Range("K132").Value = 21383.8638
$ws.Range("L132").Value = 12012.3999
$ws.Range("M132").Value = -18853.8638
$ws.Range("N132").Value = -17072.3999

$ws.Range("H133").Value = 34590.2
$ws.Range("J133").Value = 34590.2
$ws.Range("L133").Value = 34590.2
$ws.Range("N133").Value = -39650.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2006390.1
$ws.Range("I5").Value = 20000000
$ws.Range("J5").Value = 7100.1113
$ws.Range("K5").Value = 20000000
$ws.Range("L5").Value = 7100.1113
$ws.Range("M5").Value = -19999888
$ws.Range("N5").Value = -7324.1113

$ws.Range("H132").Value = 10003991
$ws.Range("I132").Value = 16671190
$ws.Range("J132").Value = 3192.1
$ws.Range("K132").Value = 50013570
$ws.Range("L132").Value = 9576.299999999999
$ws.Range("M132").Value = -50011040
$ws.Range("N132").Value = -14636.3
